# Y2_B2526_Biochemistry_LAB_CBL_System upload refresh:
#  - Trim the trailing space from the "Biochemistry Lab/CBL " shared string
#    (the whole Subject column uses this one string, so re-writing the full
#    column in one shot lets Excel fold it back down to a single, deduped,
#    shared-string entry without the trailing space).
#  - Move the active selection from A1:XFD1 to B4, matching the saved
#    cursor position in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Biochemistry Lab/CBL " shared string (drop trailing space) ---
# Column B (rows 2-479) all hold the same shared string; rewriting the whole
# range at once collapses them back onto a single shared-string entry.
$lastRow = 479
try {
    $found = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row
    if ($found -ge 2) { $lastRow = $found }
} catch {
    $lastRow = 479
}
$ws.Range("B2:B$lastRow").Value = "Biochemistry Lab/CBL"

# --- Move the selection to B4 ---
$ws.Range("B4").Select() | Out-Null
